$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.959.75"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.870.14"
$ws.Range("E3").Value = "  -2.37%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.12"
$ws.Range("E5").Value = "  -1.15%  "
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4992"
$ws.Range("E7").Value = "  -3.21%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3815"
$ws.Range("E8").Value = "  -4.37%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08941"
$ws.Range("E9").Value = "  -8.82%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.120"
$ws.Range("E10").Value = "  -2.80%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.41"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.314"
$ws.Range("E12").Value = "  -3.06%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.69"
$ws.Range("E13").Value = "  -2.38%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.864.30"
$ws.Range("E14").Value = "  -2.61%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.241"
$ws.Range("E15").Value = "  -2.99%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.000"
$ws.Range("E16").Value = "  -0.12%  "
$ws.Range("E17").Value = "  -3.13%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "90.84"
$ws.Range("E18").Value = "  -4.15%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06641"
$ws.Range("E19").Value = "  -0.29%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.92"
$ws.Range("E20").Value = "  -1.74%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.0000"
$ws.Range("E21").Value = "  -0.01%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.108"
$ws.Range("E22").Value = "  -3.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.977.43"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.50"
$ws.Range("E24").Value = "  -0.20%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.284"
$ws.Range("E25").Value = "  -1.11%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.076.46"
$ws.Range("E26").Value = "  -2.85%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.517"
$ws.Range("E27").Value = "  -6.43%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "158.11"
$ws.Range("E28").Value = "  +0.40%  "
$ws.Range("E29").Value = "  -2.82%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.27"
$ws.Range("E30").Value = "  -2.70%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.1059"
$ws.Range("E31").Value = "  -1.49%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.056"
$ws.Range("E32").Value = "  -5.74%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.577"
$ws.Range("E33").Value = "  -2.69%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.590"
$ws.Range("E34").Value = "  -1.21%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.365"
$ws.Range("E35").Value = "  -5.05%  "
$ws.Range("E36").Value = "  -3.54%  "
$ws.Range("E37").Value = "  -1.77%  "
$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.306"
$ws.Range("E38").Value = "  +9.77%  "
$ws.Range("B39").Value = "Algorand"
$ws.Range("C39").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2194"
$ws.Range("E39").Value = "  -1.76%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.202"
$ws.Range("E40").Value = "  -5.75%  "
$ws.Range("E41").Value = "  -1.25%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6374"
$ws.Range("E42").Value = "  -1.73%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.902"
$ws.Range("E43").Value = "  -3.61%  "
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.26"
$ws.Range("E45").Value = "  -2.09%  "
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6006"
$ws.Range("E46").Value = "  -1.54%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.279"
$ws.Range("E47").Value = "  -0.68%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.672"
$ws.Range("E48").Value = "  -2.84%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.226"
$ws.Range("E49").Value = "  +1.63%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.980"
$ws.Range("E50").Value = "  -4.39%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "121.57"
$ws.Range("E51").Value = "  -2.66%  "
